# Weekly price update: a new week's record for
# "Feria Lagunitas de Puerto Montt - Coliflor" is inserted at row 470,
# pushing the existing rows 470-506 down to 471-507.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 470 - this shifts rows 470:506
# down to 471:507 and extends the used range to A1:R507.
$ws.Rows.Item(470).Insert()

# Populate the newly inserted row 470 with the new week's values.
$ws.Cells.Item(470, 1).Value  = 4
$ws.Cells.Item(470, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(470, 3).Value  = "Los Lagos"
$ws.Cells.Item(470, 4).Value  = 45013
$ws.Cells.Item(470, 5).Value  = 10
$ws.Cells.Item(470, 6).Value  = 100112008
$ws.Cells.Item(470, 7).Value  = "Coliflor"
$ws.Cells.Item(470, 8).Value  = "Sin especificar"
$ws.Cells.Item(470, 9).Value  = "Primera"
$ws.Cells.Item(470, 10).Value = 1200
$ws.Cells.Item(470, 11).Value = 1700
$ws.Cells.Item(470, 12).Value = 1700
$ws.Cells.Item(470, 13).Value = 1700
$ws.Cells.Item(470, 14).Value = "$/unidad"
$ws.Cells.Item(470, 15).Value = "Región Metropolitana"
$ws.Cells.Item(470, 16).Value = 1700
$ws.Cells.Item(470, 17).Value = 1
$ws.Cells.Item(470, 18).Value = "Hortaliza"
